# Adds a new "Sheet1" worksheet (after "Phase 4 Generation") summarising
# quarter-turn / half-turn phase generation statistics by cube depth.

$wb = $excel.ActiveWorkbook

# --- Create the new worksheet as the last tab ------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet1"

# --- Header block (rows 3-5) ------------------------------------------------
# Values are written in the same order the original author entered them
# (this keeps the shared-strings table build order stable / diff-friendly).
$ws.Range("C3").Value = "Quarter Turns"
$ws.Range("H3").Value = "Half Turns"
$ws.Range("B3").Value = "Position"

$ws.Range("D5").Value = "Count"
$ws.Range("C5").Value = "Depth"
$ws.Range("F5").Value = "Count"
$ws.Range("E5").Value = "Depth"
$ws.Range("I5").Value = "Count"
$ws.Range("H5").Value = "Depth"
$ws.Range("K5").Value = "Count"
$ws.Range("J5").Value = "Depth"

# --- Data rows (6-10) --------------------------------------------------------
$ws.Range("B6").Value = "OBROWROGRWWWBRBWWWGOGWOYGGGWRYBBBYYYBOBYYYGRGOGROYROBR"
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 192
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 21
$ws.Range("G6").Value = 0.15

$ws.Range("B7").Value = "WBWWWWWGWOOOGWGRRRBWBBOGOGRGRBRBOOOOGYGRRRBYBYGYYYYYBY"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 30
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 237
$ws.Range("G7").Value = 0.218

$ws.Range("B8").Value = "YWRBWYOOWOOBYYOGBYBRGGOGWGWBRYGBOWRRGWBRRYGRBWOWYYBOGR"
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 1227
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 286252
$ws.Range("G8").Value = 180.633
$ws.Range("H8").Value = 2
$ws.Range("I8").Value = 246
$ws.Range("J8").Value = 6
$ws.Range("K8").Value = 507904
$ws.Range("L8").Value = 418.133

$ws.Range("B9").Value = "OGOYWYRBRWOWGOBWRWBRGGOBWGYBRGWBYYOYBRGYRYGOBOBOWYWRGR"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 237
$ws.Range("G9").Value = 0.187

$ws.Range("B10").Value = "GYWBWRBOYWRWRWRGWBOGRBORYGRGRBOBWWGBYYOYOOGYORBBWYGGOY"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 30537
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 21814238
$ws.Range("G10").Value = 13324.681

# --- Row 4 phase labels (added last, as in the original edit sequence) ------
$ws.Range("C4").Value = "Phase 1"
$ws.Range("E4").Value = "Phase 2"
$ws.Range("G4").Value = "Time"
$ws.Range("H4").Value = "Phase 1"
$ws.Range("J4").Value = "Phase 2"
$ws.Range("L4").Value = "Time"

# --- Merge the header cells --------------------------------------------------
$ws.Range("B3:B5").Merge()
$ws.Range("C3:G3").Merge()
$ws.Range("H3:L3").Merge()
$ws.Range("C4:D4").Merge()
$ws.Range("E4:F4").Merge()
$ws.Range("G4:G5").Merge()
$ws.Range("H4:I4").Merge()
$ws.Range("J4:K4").Merge()
$ws.Range("L4:L5").Merge()

# --- Centre-align the whole header block (rows 3-5) -------------------------
$headerRange = $ws.Range("B3:L5")
$headerRange.VerticalAlignment = -4108
$headerRange.HorizontalAlignment = -4108

# --- Column B sized to fit the 54-character cube-position strings -----------
$ws.Columns.Item(2).ColumnWidth = 73

# --- Match the author's final selection/view state --------------------------
$ws.Range("F9").Select()
